# Update the "within 100" arithmetic drill table: each of the 100 cells
# (20 rows x 5 columns) gets its equation text replaced, cell by cell, in
# row-major order. Setting Cell.Range.Text preserves the cell's existing
# run formatting (font/size) and the trailing cell-end mark.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "48-28="
$t.Cell(1, 2).Range.Text = "95-80="
$t.Cell(1, 3).Range.Text = "86+1="
$t.Cell(1, 4).Range.Text = "2+59="
$t.Cell(1, 5).Range.Text = "41+56="
$t.Cell(2, 1).Range.Text = "66-58="
$t.Cell(2, 2).Range.Text = "18+25="
$t.Cell(2, 3).Range.Text = "84-23="
$t.Cell(2, 4).Range.Text = "27+47="
$t.Cell(2, 5).Range.Text = "90-6="
$t.Cell(3, 1).Range.Text = "72-32="
$t.Cell(3, 2).Range.Text = "15-12="
$t.Cell(3, 3).Range.Text = "60-13="
$t.Cell(3, 4).Range.Text = "55-18="
$t.Cell(3, 5).Range.Text = "67+14="
$t.Cell(4, 1).Range.Text = "27+2="
$t.Cell(4, 2).Range.Text = "82-44="
$t.Cell(4, 3).Range.Text = "10+3="
$t.Cell(4, 4).Range.Text = "31+31="
$t.Cell(4, 5).Range.Text = "22+21="
$t.Cell(5, 1).Range.Text = "31+31="
$t.Cell(5, 2).Range.Text = "27-16="
$t.Cell(5, 3).Range.Text = "10-7="
$t.Cell(5, 4).Range.Text = "46+17="
$t.Cell(5, 5).Range.Text = "76-68="
$t.Cell(6, 1).Range.Text = "37+1="
$t.Cell(6, 2).Range.Text = "77-11="
$t.Cell(6, 3).Range.Text = "69-63="
$t.Cell(6, 4).Range.Text = "55-34="
$t.Cell(6, 5).Range.Text = "99-62="
$t.Cell(7, 1).Range.Text = "67-24="
$t.Cell(7, 2).Range.Text = "85-24="
$t.Cell(7, 3).Range.Text = "47+21="
$t.Cell(7, 4).Range.Text = "62-41="
$t.Cell(7, 5).Range.Text = "22+61="
$t.Cell(8, 1).Range.Text = "31-8="
$t.Cell(8, 2).Range.Text = "20-11="
$t.Cell(8, 3).Range.Text = "42+19="
$t.Cell(8, 4).Range.Text = "44+6="
$t.Cell(8, 5).Range.Text = "94-65="
$t.Cell(9, 1).Range.Text = "64-50="
$t.Cell(9, 2).Range.Text = "70-41="
$t.Cell(9, 3).Range.Text = "49+8="
$t.Cell(9, 4).Range.Text = "78+1="
$t.Cell(9, 5).Range.Text = "29+17="
$t.Cell(10, 1).Range.Text = "35-14="
$t.Cell(10, 2).Range.Text = "20+1="
$t.Cell(10, 3).Range.Text = "42+10="
$t.Cell(10, 4).Range.Text = "73-9="
$t.Cell(10, 5).Range.Text = "85-73="
$t.Cell(11, 1).Range.Text = "41+24="
$t.Cell(11, 2).Range.Text = "1+63="
$t.Cell(11, 3).Range.Text = "26+45="
$t.Cell(11, 4).Range.Text = "31+32="
$t.Cell(11, 5).Range.Text = "6+78="
$t.Cell(12, 1).Range.Text = "75-68="
$t.Cell(12, 2).Range.Text = "4+4="
$t.Cell(12, 3).Range.Text = "5-4="
$t.Cell(12, 4).Range.Text = "35+41="
$t.Cell(12, 5).Range.Text = "37-21="
$t.Cell(13, 1).Range.Text = "1+93="
$t.Cell(13, 2).Range.Text = "76-15="
$t.Cell(13, 3).Range.Text = "84-67="
$t.Cell(13, 4).Range.Text = "96-76="
$t.Cell(13, 5).Range.Text = "74+15="
$t.Cell(14, 1).Range.Text = "1+37="
$t.Cell(14, 2).Range.Text = "17+6="
$t.Cell(14, 3).Range.Text = "67-55="
$t.Cell(14, 4).Range.Text = "95-33="
$t.Cell(14, 5).Range.Text = "23+33="
$t.Cell(15, 1).Range.Text = "70+19="
$t.Cell(15, 2).Range.Text = "85-57="
$t.Cell(15, 3).Range.Text = "4+73="
$t.Cell(15, 4).Range.Text = "62-54="
$t.Cell(15, 5).Range.Text = "95-11="
$t.Cell(16, 1).Range.Text = "52+31="
$t.Cell(16, 2).Range.Text = "95-17="
$t.Cell(16, 3).Range.Text = "35+12="
$t.Cell(16, 4).Range.Text = "16+59="
$t.Cell(16, 5).Range.Text = "36-3="
$t.Cell(17, 1).Range.Text = "92-79="
$t.Cell(17, 2).Range.Text = "27+15="
$t.Cell(17, 3).Range.Text = "21+32="
$t.Cell(17, 4).Range.Text = "51-33="
$t.Cell(17, 5).Range.Text = "75-6="
$t.Cell(18, 1).Range.Text = "94-53="
$t.Cell(18, 2).Range.Text = "47-29="
$t.Cell(18, 3).Range.Text = "18+8="
$t.Cell(18, 4).Range.Text = "44-32="
$t.Cell(18, 5).Range.Text = "56-5="
$t.Cell(19, 1).Range.Text = "60+22="
$t.Cell(19, 2).Range.Text = "36+1="
$t.Cell(19, 3).Range.Text = "54-16="
$t.Cell(19, 4).Range.Text = "70-15="
$t.Cell(19, 5).Range.Text = "56-30="
$t.Cell(20, 1).Range.Text = "21+32="
$t.Cell(20, 2).Range.Text = "4+42="
$t.Cell(20, 3).Range.Text = "77-53="
$t.Cell(20, 4).Range.Text = "83-31="
$t.Cell(20, 5).Range.Text = "14+79="
